$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.880.30'
$ws.Range("E2").Value = '  -1.70%  '
$ws.Range("D3").Value = '2.603.83'
$ws.Range("E3").Value = '  -0.93%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.36%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.553'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.61%  '
$ws.Range("D9").Value = '2.602.78'
$ws.Range("E9").Value = '  -1.01%  '
$ws.Range("E10").Value = '  -2.65%  '
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("E12").Value = '  -2.15%  '
$ws.Range("E13").Value = '  -4.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.86%  '
$ws.Range("D15").Value = '3.078.71'
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("E16").Value = '  -4.97%  '
$ws.Range("D17").Value = '66.760.32'
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").Value = '2.599.35'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '363.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.11%  '
$ws.Range("E21").Value = '  -6.13%  '
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("E25").Value = '  -2.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '67.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.17%  '
$ws.Range("E27").Value = '  -0.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '576.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.67%  '
$ws.Range("D30").Value = '0.0₃0998'
$ws.Range("E30").Value = '  -5.47%  '
$ws.Range("E31").Value = '  -6.50%  '
$ws.Range("E32").Value = '  -4.62%  '
$ws.Range("E33").Value = '  -3.02%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.122'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.49%  '
$ws.Range("E36").Value = '  -4.54%  '
$ws.Range("E37").Value = '  -2.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '156.09'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.78%  '
$ws.Range("E40").Value = '  -2.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.21'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.83%  '
$ws.Range("E42").Value = '  -4.42%  '
$ws.Range("E43").Value = '  -4.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.94'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.35'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '154.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.63%  '
$ws.Range("E48").Value = '  -3.73%  '
$ws.Range("E49").Value = '  -1.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.45'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.49%  '
$ws.Range("E51").Value = '  -3.27%  '
